$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15, shifting existing rows 15-76 down to 16-77.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 44804
$ws.Range("D15").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100112029
$ws.Range("G15").Value = "Orégano"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 18000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 18000
$ws.Range("N15").Value = "$/docena de atados"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 6000
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = "Hortaliza"
